$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the path stored in D2 (RutaPagaduria value) to the new relative path
$ws.Range("D2").Value = '"src/test/resources/Data/PagaduriaAplicacion/"'

# Move the active selection from F2 to D4, matching the saved view state
$ws.Range("D4").Select()
